$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.01028766666666667
$ws.Range("H2").Value = 0.030863
$ws.Range("I2").Value = 0.004475330795722314
$ws.Range("J2").Value = 0.004475330795722313
$ws.Range("M2").Value = 0.013559
$ws.Range("N2").Value = 0.040677
$ws.Range("O2").Value = 0.01041720286693961
$ws.Range("P2").Value = 0.01041720286693961
$ws.Range("Q2").Value = 0.0001394904723333334
$ws.Range("R2").Value = 0.001255414251
$ws.Range("S2").Value = 0.00004662042879570163
$ws.Range("T2").Value = 0.00004662042879570163

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.01028766666666667
$ws.Range("H3").Value = 0.030863
$ws.Range("I3").Value = 0.004475330795722314
$ws.Range("J3").Value = 0.004475330795722313
$ws.Range("M3").Value = 1.288038
$ws.Range("N3").Value = 3.864114
$ws.Range("O3").Value = 0.9895827971330603
$ws.Range("P3").Value = 0.9895827971330604
$ws.Range("Q3").Value = 0.013250905598
$ws.Range("R3").Value = 0.119258150382
$ws.Range("S3").Value = 0.004428710366926612
$ws.Range("T3").Value = 0.004428710366926611

# Row 4
$ws.Range("G4").Value = 2.288462333333333
$ws.Range("H4").Value = 6.865386999999999
$ws.Range("I4").Value = 0.9955246692042776
$ws.Range("J4").Value = 0.9955246692042776
$ws.Range("M4").Value = 0.013559
$ws.Range("N4").Value = 0.040677
$ws.Range("O4").Value = 0.01041720286693961
$ws.Range("P4").Value = 0.01041720286693961
$ws.Range("Q4").Value = 0.03102926077766666
$ws.Range("R4").Value = 0.279263346999
$ws.Range("S4").Value = 0.01037058243814391
$ws.Range("T4").Value = 0.01037058243814391

# Row 5
$ws.Range("G5").Value = 2.288462333333333
$ws.Range("H5").Value = 6.865386999999999
$ws.Range("I5").Value = 0.9955246692042776
$ws.Range("J5").Value = 0.9955246692042776
$ws.Range("M5").Value = 1.288038
$ws.Range("N5").Value = 3.864114
$ws.Range("O5").Value = 0.9895827971330603
$ws.Range("P5").Value = 0.9895827971330604
$ws.Range("Q5").Value = 2.947626446901999
$ws.Range("R5").Value = 26.528638022118
$ws.Range("S5").Value = 0.9851540867661336
$ws.Range("T5").Value = 0.9851540867661337
